$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6th column).
# This shifts old F -> G, old G -> H, carrying values/styles along.
$ws.Columns("F").Insert()

# New column F should have the same display width as column E.
$ws.Range("F1").EntireColumn.ColumnWidth = $ws.Range("E1").EntireColumn.ColumnWidth

# New header cell F1 = "Staff Req"
$ws.Range("F1").Value = "Staff Req"

# New formula column F2:F7 = G + H (Kitchen Req + Service Req)
$ws.Range("F2").Formula = "=G2+H2"
$ws.Range("F3").Formula = "=G3+H3"
$ws.Range("F4").Formula = "=G4+H4"
$ws.Range("F5").Formula = "=G5+H5"
$ws.Range("F6").Formula = "=G6+H6"
$ws.Range("F7").Formula = "=G7+H7"

# Match style/format of other numeric data cells (style index 2 in the target).
$ws.Range("A2").Copy()
$ws.Range("F2:F7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to F3 to match the saved view state.
$ws.Range("F3").Select()
